$wb = $excel.ActiveWorkbook

# Add the new "properties" worksheet
$propsSheet = $wb.Worksheets.Add()
$propsSheet.Name = "properties"

$headers = @("partition", "aspect", "key", "type", "value")
for ($c = 1; $c -le $headers.Length; $c++) {
    $propsSheet.Cells.Item(1, $c).Value = $headers[$c - 1]
}

$rows = @(
    @("Table", "default", "colOrder", "array", '["plot_name","location_latitude","location_longitude","location_altitude","location_accuracy","planting"]'),
    @("Table", "default", "defaultViewType", "string", "MAP"),
    @("Table", "default", "detailViewFileName", "configpath", "config/tables/plot/html/plot_detail.html"),
    @("Table", "default", "listViewFileName", "configpath", "config/tables/plot/html/plot_list.html"),
    @("Table", "default", "mapListViewFileName", "configpath", "config/tables/plot/html/plot_list.html"),
    @("TableMapFragment", "default", "keyColorRuleType", "string", "None"),
    @("TableMapFragment", "default", "keyMapLatCol", "string", "location_latitude"),
    @("TableMapFragment", "default", "keyMapLongCol", "string", "location_longitude")
)

$r = 2
foreach ($row in $rows) {
    for ($c = 1; $c -le $row.Length; $c++) {
        $propsSheet.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}


# Move the new sheet to the end (after "settings"), now that all data is written
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$propsSheet.Move($null, $lastSheet) | Out-Null

$finalPropsSheet = $wb.Worksheets.Item("properties")
$finalPropsSheet.Activate() | Out-Null
$finalPropsSheet.Rows.Item(3).Select() | Out-Null
